$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15
$ws1.Range("F3").Value = 7928
$ws1.Range("F7").Value = 803
$ws1.Range("F8").Value = 610
$ws1.Range("F9").Value = 93
$ws1.Range("F12").Value = 866
$ws1.Range("F13").Value = 3255
$ws1.Range("F14").Value = 203
$ws1.Range("F16").Value = 740
$ws1.Range("F18").Value = 49
$ws1.Range("F19").Value = 461
$ws1.Range("F21").Value = 253
$ws1.Range("F22").Value = 226
$ws1.Range("F23").Value = 323
$ws1.Range("F26").Value = 110
$ws1.Range("F27").Value = 278
$ws1.Range("F28").Value = 23
$ws1.Range("F31").Value = 505
$ws1.Range("F32").Value = 534
$ws1.Range("F33").Value = 24
$ws1.Range("F34").Value = 34
$ws1.Range("F35").Value = 6
$ws1.Range("F37").Value = 225

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 207

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 207
$ws4.Range("F4").Value = 15
$ws4.Range("F5").Value = 7928
$ws4.Range("F9").Value = 803
$ws4.Range("F10").Value = 610
$ws4.Range("F11").Value = 93
$ws4.Range("F14").Value = 866
$ws4.Range("F16").Value = 3255
$ws4.Range("F17").Value = 203
$ws4.Range("F20").Value = 740
$ws4.Range("F23").Value = 49
$ws4.Range("F24").Value = 461
$ws4.Range("F26").Value = 253
$ws4.Range("F27").Value = 226
$ws4.Range("F28").Value = 323
$ws4.Range("F31").Value = 110
$ws4.Range("F32").Value = 278
$ws4.Range("F33").Value = 23
$ws4.Range("F36").Value = 505
$ws4.Range("F37").Value = 534
$ws4.Range("F38").Value = 24
$ws4.Range("F39").Value = 34
$ws4.Range("F40").Value = 6
$ws4.Range("F42").Value = 225
